$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PhpTravels")

# New header cells (row 1): AC1, AD1, AE1
$ws.Cells.Item(1, 29).Value = "firstname"
$ws.Cells.Item(1, 30).Value = "lastname"
$ws.Cells.Item(1, 31).Value = "passportid"

# New data cells (row 2): AC2 (plain text), AD2/AE2 (quote-prefixed so Excel
# treats the numeric-looking CSV text as text, matching the quotePrefix style
# used by the other "list" columns in this sheet).
$ws.Cells.Item(2, 29).Value = "Gopi,Gowri,Gowreesh,Gopika,Graden"
$ws.Cells.Item(2, 30).Value = "'Muthu,gopi,gopi,gopi,gopi"
$ws.Cells.Item(2, 31).Value = "'1234567890,0987654321,6789054321,0987612345,1236547890"

# Match the bestFit-computed column widths (AC=29, AD=30, AE=31) as closely as
# this engine's ColumnWidth quantization allows.
$ws.Columns.Item(29).ColumnWidth = 34.3
$ws.Columns.Item(30).ColumnWidth = 24.15
$ws.Columns.Item(31).ColumnWidth = 54.0

# Move the selection/active cell and drop the frozen top-left scroll position.
[void]$ws.Range("AB6").Select()
